# Update to recent changes in local error codes; switch to testing preflight
# with ErrorCodes.xlsx

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Errors_")

# ---------------------------------------------------------------------------
# 1. Update the error-code values in column A (rows 10-26 were renumbered,
#    rows 27-32 were rebased from 2xx/3xx to 4xx/5xx).
# ---------------------------------------------------------------------------
$ws.Range("A10").Value = 121
$ws.Range("A11").Value = 130
$ws.Range("A12").Value = 136
$ws.Range("A13").Value = 150
$ws.Range("A14").Value = 155
$ws.Range("A15").Value = 170
$ws.Range("A16").Value = 174
$ws.Range("A17").Value = 190
$ws.Range("A18").Value = 193
$ws.Range("A19").Value = 210
$ws.Range("A20").Value = 212
$ws.Range("A21").Value = 230
$ws.Range("A22").Value = 237
$ws.Range("A23").Value = 250
$ws.Range("A24").Value = 258
$ws.Range("A25").Value = 270
$ws.Range("A26").Value = 271

$ws.Range("A27").Value = 400
$ws.Range("A28").Value = 401
$ws.Range("A29").Value = 500
$ws.Range("A30").Value = 501
$ws.Range("A31").Value = 510
$ws.Range("A32").Value = 511

# ---------------------------------------------------------------------------
# 2. Re-style the "ColAllNumeric" through "NoDuplicateCols" blocks (rows
#    13-26) with the 60% - Accent5 cell style.
# ---------------------------------------------------------------------------
$ws.Range("A13:A26").Style = "60% - Accent5"
$ws.Range("A13:A26").Font.Color = 0

# ---------------------------------------------------------------------------
# 3. Tidy up the error message text - it no longer calls out "col_b" by name.
# ---------------------------------------------------------------------------
$ws.Range("D14").Value = "ERROR: Column must contain only non-blank numeric values"

# ---------------------------------------------------------------------------
# 4. Move the active selection to A11 (cursor position when the file was
#    last saved).
# ---------------------------------------------------------------------------
$ws.Range("A11").Select()
